$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a new paragraph right after "Assignment 04" that contains only a
#    hyperlink to the author's GitHub profile (mirrors the paragraph that used
#    to live at the end of the document).
# ---------------------------------------------------------------------------
$assignmentPara = $d.Paragraphs.Item(4)

# InsertParagraphAfter on the paragraph's own Range (not a collapsed point at
# its end) keeps the new paragraph on the "Normal" style instead of picking
# up the following paragraph's style.
$assignmentPara.Range.InsertParagraphAfter()

$githubPara = $d.Paragraphs.Item($assignmentPara.Index + 1)

$d.Hyperlinks.Add($githubPara.Range, "https://github.com/jordan-aloysius", `
    [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, `
    "https://github.com/jordan-aloysius") | Out-Null

# ---------------------------------------------------------------------------
# 2. Remove the old "Link to my GitHub: ..." paragraph at the end of the
#    document (text run + the two hyperlink runs), merging what is left back
#    into the preceding "Summary" paragraph so no empty paragraph remains.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$summaryPara = $d.Paragraphs.Item($count - 1)
$linkPara = $d.Paragraphs.Item($count)

$boundary = $summaryPara.Range.End - 1

# Delete just the paragraph mark that separates the two paragraphs so they
# merge into a single paragraph (keeping the formatting of the first one).
$d.Range($summaryPara.Range.End - 1, $summaryPara.Range.End).Delete()

$mergedPara = $d.Paragraphs.Item($d.Paragraphs.Count)
# Delete everything that used to be the "Link to my GitHub" paragraph's
# content, but keep the final paragraph mark of the document.
$d.Range($boundary, $mergedPara.Range.End - 1).Delete()

# ---------------------------------------------------------------------------
# 3. Flip "Allow punctuation to extend past text extents" (w:overflowPunct)
#    on for both the Normal and No Spacing styles.
# ---------------------------------------------------------------------------
$d.Styles.Item("Normal").ParagraphFormat.HangingPunctuation = $true
$d.Styles.Item("NoSpacing").ParagraphFormat.HangingPunctuation = $true

Write-Output "done"
